# EPS v3.3.1 -> v3.4.2 update
# - Coal (hard coal) is no longer bid at peak capacity factor for the US: set to 0.
# - Add explanatory notes about this change to the "About" sheet.
# - Update selection on the data sheet.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("BDSBaPCF")

# --- Data sheet: set "hard coal" bid-at-peak flag to 0 (B13 recalculates via formula =B2) ---
$wsData.Range("B2").Value = 0

# --- About sheet: append the new explanatory note rows ---
$wsAbout.Range("A24").Value = "For the United States, we have set coal to 0 as of version 3.4. This reflects"
$wsAbout.Range("A25").Value = "the fact that certain air quality / environmental restrictions, as well as current"
$wsAbout.Range("A26").Value = "supply chain logistics, limit the amount the coal dispatches annually. "

# --- Update selections to match the saved workbook state ---
$wsData.Range("B4").Select()
$wsAbout.Range("A27").Select()
$wsAbout.Activate()
